$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 18.60499954223633
$ws.Cells.Item(2, 5).Value = 21.0890007019043
$ws.Cells.Item(2, 6).Value = 22.63249969482422
$ws.Cells.Item(2, 7).Value = 18.41699981689453
$ws.Cells.Item(2, 8).Value = 10664912097
$ws.Cells.Item(2, 9).Value = "AMZN"

$ws.Cells.Item(3, 4).Value = 21.96750068664551
$ws.Cells.Item(3, 5).Value = 26.8075008392334
$ws.Cells.Item(3, 6).Value = 29.02849960327148
$ws.Cells.Item(3, 7).Value = 21.27849960327148
$ws.Cells.Item(3, 8).Value = 10664912097
$ws.Cells.Item(3, 9).Value = "AMZN"

$ws.Cells.Item(4, 4).Value = 25.54999923706055
$ws.Cells.Item(4, 5).Value = 31.29500007629395
$ws.Cells.Item(4, 6).Value = 31.5359992980957
$ws.Cells.Item(4, 7).Value = 25.29999923706055
$ws.Cells.Item(4, 8).Value = 10664912097
$ws.Cells.Item(4, 9).Value = "AMZN"

$ws.Cells.Item(5, 4).Value = 32.81449890136719
$ws.Cells.Item(5, 5).Value = 29.35000038146973
$ws.Cells.Item(5, 6).Value = 32.88600158691406
$ws.Cells.Item(5, 7).Value = 27.35899925231934
$ws.Cells.Item(5, 8).Value = 10664912097
$ws.Cells.Item(5, 9).Value = "AMZN"

$ws.Cells.Item(6, 4).Value = 29.52449989318848
$ws.Cells.Item(6, 5).Value = 32.97949981689453
$ws.Cells.Item(6, 6).Value = 33.49900054931641
$ws.Cells.Item(6, 7).Value = 29.26250076293945
$ws.Cells.Item(6, 8).Value = 10664912097
$ws.Cells.Item(6, 9).Value = "AMZN"

$ws.Cells.Item(7, 4).Value = 35.86600112915039
$ws.Cells.Item(7, 5).Value = 37.94049835205078
$ws.Cells.Item(7, 6).Value = 38.29999923706055
$ws.Cells.Item(7, 7).Value = 35.82699966430664
$ws.Cells.Item(7, 8).Value = 10664912097
$ws.Cells.Item(7, 9).Value = "AMZN"

$ws.Cells.Item(8, 4).Value = 41.79999923706055
$ws.Cells.Item(8, 5).Value = 39.49100112915039
$ws.Cells.Item(8, 6).Value = 42.36050033569336
$ws.Cells.Item(8, 7).Value = 38.73049926757812
$ws.Cells.Item(8, 8).Value = 10664912097
$ws.Cells.Item(8, 9).Value = "AMZN"

$ws.Cells.Item(9, 4).Value = 37.89599990844727
$ws.Cells.Item(9, 5).Value = 41.17399978637695
$ws.Cells.Item(9, 6).Value = 42.19200134277344
$ws.Cells.Item(9, 7).Value = 37.3849983215332
$ws.Cells.Item(9, 8).Value = 10664912097
$ws.Cells.Item(9, 9).Value = "AMZN"

$ws.Cells.Item(10, 4).Value = 44.40000152587891
$ws.Cells.Item(10, 5).Value = 46.2495002746582
$ws.Cells.Item(10, 6).Value = 47.47949981689453
$ws.Cells.Item(10, 7).Value = 44.2244987487793
$ws.Cells.Item(10, 8).Value = 10664912097
$ws.Cells.Item(10, 9).Value = "AMZN"

$ws.Cells.Item(11, 4).Value = 48.63949966430664
$ws.Cells.Item(11, 5).Value = 49.38899993896485
$ws.Cells.Item(11, 6).Value = 54.16550064086914
$ws.Cells.Item(11, 7).Value = 47.54999923706055
$ws.Cells.Item(11, 8).Value = 10664912097
$ws.Cells.Item(11, 9).Value = "AMZN"

$ws.Cells.Item(12, 4).Value = 48.20000076293945
$ws.Cells.Item(12, 5).Value = 55.26399993896485
$ws.Cells.Item(12, 6).Value = 56.13949966430664
$ws.Cells.Item(12, 7).Value = 47.51850128173828
$ws.Cells.Item(12, 8).Value = 10664912097
$ws.Cells.Item(12, 9).Value = "AMZN"

$ws.Cells.Item(13, 4).Value = 58.59999847412109
$ws.Cells.Item(13, 5).Value = 72.54450225830078
$ws.Cells.Item(13, 6).Value = 73.62899780273438
$ws.Cells.Item(13, 7).Value = 58.5255012512207
$ws.Cells.Item(13, 8).Value = 10664912097
$ws.Cells.Item(13, 9).Value = "AMZN"

$ws.Cells.Item(14, 4).Value = 70.88099670410156
$ws.Cells.Item(14, 5).Value = 78.30650329589844
$ws.Cells.Item(14, 6).Value = 81.90499877929688
$ws.Cells.Item(14, 7).Value = 67.64399719238281
$ws.Cells.Item(14, 8).Value = 10664912097
$ws.Cells.Item(14, 9).Value = "AMZN"

$ws.Cells.Item(15, 4).Value = 84.13500213623047
$ws.Cells.Item(15, 5).Value = 88.87200164794922
$ws.Cells.Item(15, 6).Value = 94.00250244140624
$ws.Cells.Item(15, 7).Value = 83.90299987792969
$ws.Cells.Item(15, 8).Value = 10664912097
$ws.Cells.Item(15, 9).Value = "AMZN"

$ws.Cells.Item(16, 4).Value = 101.0995025634766
$ws.Cells.Item(16, 5).Value = 79.90049743652344
$ws.Cells.Item(16, 6).Value = 101.6595001220703
$ws.Cells.Item(16, 7).Value = 73.81800079345703
$ws.Cells.Item(16, 8).Value = 10664912097
$ws.Cells.Item(16, 9).Value = "AMZN"

$ws.Cells.Item(17, 4).Value = 73.26000213623047
$ws.Cells.Item(17, 5).Value = 85.93650054931641
$ws.Cells.Item(17, 6).Value = 86.82050323486328
$ws.Cells.Item(17, 7).Value = 73.04650115966797
$ws.Cells.Item(17, 8).Value = 10664912097
$ws.Cells.Item(17, 9).Value = "AMZN"

$ws.Cells.Item(18, 4).Value = 90.00550079345703
$ws.Cells.Item(18, 5).Value = 96.32599639892578
$ws.Cells.Item(18, 6).Value = 97.81700134277344
$ws.Cells.Item(18, 7).Value = 89.93650054931641
$ws.Cells.Item(18, 8).Value = 10664912097
$ws.Cells.Item(18, 9).Value = "AMZN"

$ws.Cells.Item(19, 4).Value = 96.14900207519533
$ws.Cells.Item(19, 5).Value = 93.33899688720705
$ws.Cells.Item(19, 6).Value = 101.7900009155273
$ws.Cells.Item(19, 7).Value = 92.47200012207033
$ws.Cells.Item(19, 8).Value = 10664912097
$ws.Cells.Item(19, 9).Value = "AMZN"

$ws.Cells.Item(20, 4).Value = 87.30000305175781
$ws.Cells.Item(20, 5).Value = 88.83300018310547
$ws.Cells.Item(20, 6).Value = 89.94249725341797
$ws.Cells.Item(20, 7).Value = 84.25299835205078
$ws.Cells.Item(20, 8).Value = 10664912097
$ws.Cells.Item(20, 9).Value = "AMZN"

$ws.Cells.Item(21, 4).Value = 93.75
$ws.Cells.Item(21, 5).Value = 100.4359970092773
$ws.Cells.Item(21, 6).Value = 102.786003112793
$ws.Cells.Item(21, 7).Value = 90.76699829101562
$ws.Cells.Item(21, 8).Value = 10664912097
$ws.Cells.Item(21, 9).Value = "AMZN"

$ws.Cells.Item(22, 4).Value = 96.64849853515624
$ws.Cells.Item(22, 5).Value = 123.6999969482422
$ws.Cells.Item(22, 6).Value = 123.75
$ws.Cells.Item(22, 7).Value = 94.45749664306641
$ws.Cells.Item(22, 8).Value = 10664912097
$ws.Cells.Item(22, 9).Value = "AMZN"

$ws.Cells.Item(23, 4).Value = 137.8995056152344
$ws.Cells.Item(23, 5).Value = 158.2339935302734
$ws.Cells.Item(23, 6).Value = 167.2144927978516
$ws.Cells.Item(23, 7).Value = 137.6999969482422
$ws.Cells.Item(23, 8).Value = 10664912097
$ws.Cells.Item(23, 9).Value = "AMZN"

$ws.Cells.Item(24, 4).Value = 160.3999938964844
$ws.Cells.Item(24, 5).Value = 151.8074951171875
$ws.Cells.Item(24, 6).Value = 174.8119964599609
$ws.Cells.Item(24, 7).Value = 150.9499969482422
$ws.Cells.Item(24, 8).Value = 10664912097
$ws.Cells.Item(24, 9).Value = "AMZN"

$ws.Cells.Item(25, 4).Value = 163.5
$ws.Cells.Item(25, 5).Value = 160.3099975585938
$ws.Cells.Item(25, 6).Value = 168.1945037841797
$ws.Cells.Item(25, 7).Value = 154.3000030517578
$ws.Cells.Item(25, 8).Value = 10664912097
$ws.Cells.Item(25, 9).Value = "AMZN"

$ws.Cells.Item(26, 4).Value = 155.8970031738281
$ws.Cells.Item(26, 5).Value = 173.3710021972656
$ws.Cells.Item(26, 6).Value = 177.6999969482422
$ws.Cells.Item(26, 7).Value = 155.7774963378906
$ws.Cells.Item(26, 8).Value = 10664912097
$ws.Cells.Item(26, 9).Value = "AMZN"

$ws.Cells.Item(27, 4).Value = 171.7304992675781
$ws.Cells.Item(27, 5).Value = 166.3795013427734
$ws.Cells.Item(27, 6).Value = 188.6540069580078
$ws.Cells.Item(27, 7).Value = 165.3489990234375
$ws.Cells.Item(27, 8).Value = 10664912097
$ws.Cells.Item(27, 9).Value = "AMZN"

$ws.Cells.Item(28, 4).Value = 164.4505004882812
$ws.Cells.Item(28, 5).Value = 168.6215057373047
$ws.Cells.Item(28, 6).Value = 173.9499969482422
$ws.Cells.Item(28, 7).Value = 158.8125
$ws.Cells.Item(28, 8).Value = 10664912097
$ws.Cells.Item(28, 9).Value = "AMZN"

$ws.Cells.Item(29, 4).Value = 167.5500030517578
$ws.Cells.Item(29, 5).Value = 149.5735015869141
$ws.Cells.Item(29, 6).Value = 171.3999938964844
$ws.Cells.Item(29, 7).Value = 135.3520050048828
$ws.Cells.Item(29, 8).Value = 10664912097
$ws.Cells.Item(29, 9).Value = "AMZN"

$ws.Cells.Item(30, 4).Value = 164.1495056152344
$ws.Cells.Item(30, 5).Value = 124.2815017700195
$ws.Cells.Item(30, 6).Value = 168.3945007324219
$ws.Cells.Item(30, 7).Value = 121.625
$ws.Cells.Item(30, 8).Value = 10664912097
$ws.Cells.Item(30, 9).Value = "AMZN"

$ws.Cells.Item(31, 4).Value = 106.2900009155273
$ws.Cells.Item(31, 5).Value = 134.9499969482422
$ws.Cells.Item(31, 6).Value = 137.6499938964844
$ws.Cells.Item(31, 7).Value = 105.8499984741211
$ws.Cells.Item(31, 8).Value = 10664912097
$ws.Cells.Item(31, 9).Value = "AMZN"

$ws.Cells.Item(32, 4).Value = 113.5800018310547
$ws.Cells.Item(32, 5).Value = 102.4400024414062
$ws.Cells.Item(32, 6).Value = 123
$ws.Cells.Item(32, 7).Value = 97.66000366210938
$ws.Cells.Item(32, 8).Value = 10664912097
$ws.Cells.Item(32, 9).Value = "AMZN"

$ws.Cells.Item(33, 4).Value = 85.45999908447266
$ws.Cells.Item(33, 5).Value = 103.129997253418
$ws.Cells.Item(33, 6).Value = 103.4899978637695
$ws.Cells.Item(33, 7).Value = 81.43000030517578
$ws.Cells.Item(33, 8).Value = 10664912097
$ws.Cells.Item(33, 9).Value = "AMZN"

$ws.Cells.Item(34, 4).Value = 102.3000030517578
$ws.Cells.Item(34, 5).Value = 105.4499969482422
$ws.Cells.Item(34, 6).Value = 110.8600006103516
$ws.Cells.Item(34, 7).Value = 97.70999908447266
$ws.Cells.Item(34, 8).Value = 10664912097
$ws.Cells.Item(34, 9).Value = "AMZN"

$ws.Cells.Item(35, 4).Value = 130.8200073242188
$ws.Cells.Item(35, 5).Value = 133.6799926757812
$ws.Cells.Item(35, 6).Value = 136.6499938964844
$ws.Cells.Item(35, 7).Value = 125.9199981689453
$ws.Cells.Item(35, 8).Value = 10664912097
$ws.Cells.Item(35, 9).Value = "AMZN"

$ws.Cells.Item(36, 4).Value = 127.2799987792969
$ws.Cells.Item(36, 5).Value = 133.0899963378906
$ws.Cells.Item(36, 6).Value = 134.4799957275391
$ws.Cells.Item(36, 7).Value = 118.3499984741211
$ws.Cells.Item(36, 8).Value = 10664912097
$ws.Cells.Item(36, 9).Value = "AMZN"

$ws.Cells.Item(37, 4).Value = 151.5399932861328
$ws.Cells.Item(37, 5).Value = 155.1999969482422
$ws.Cells.Item(37, 6).Value = 161.7299957275391
$ws.Cells.Item(37, 7).Value = 144.0500030517578
$ws.Cells.Item(37, 8).Value = 10664912097
$ws.Cells.Item(37, 9).Value = "AMZN"

$ws.Cells.Item(38, 4).Value = 180.7899932861328
$ws.Cells.Item(38, 5).Value = 175
$ws.Cells.Item(38, 6).Value = 189.7700042724609
$ws.Cells.Item(38, 7).Value = 166.3200073242188
$ws.Cells.Item(38, 8).Value = 10664912097
$ws.Cells.Item(38, 9).Value = "AMZN"

$ws.Cells.Item(39, 4).Value = 193.4900054931641
$ws.Cells.Item(39, 5).Value = 186.979995727539
$ws.Cells.Item(39, 6).Value = 201.1999969482422
$ws.Cells.Item(39, 7).Value = 176.8000030517578
$ws.Cells.Item(39, 8).Value = 10664912097
$ws.Cells.Item(39, 9).Value = "AMZN"

$ws.Cells.Item(40, 4).Value = 184.8999938964844
$ws.Cells.Item(40, 5).Value = 186.3999938964844
$ws.Cells.Item(40, 6).Value = 195.6100006103516
$ws.Cells.Item(40, 7).Value = 180.25
$ws.Cells.Item(40, 8).Value = 10664912097
$ws.Cells.Item(40, 9).Value = "AMZN"

$ws.Cells.Item(41, 4).Value = 222.0299987792969
$ws.Cells.Item(41, 5).Value = 237.6799926757812
$ws.Cells.Item(41, 6).Value = 241.7700042724609
$ws.Cells.Item(41, 7).Value = 216.1999969482422
$ws.Cells.Item(41, 8).Value = 10664912097
$ws.Cells.Item(41, 9).Value = "AMZN"

$ws.Cells.Item(43, 4).Value = 219.5
$ws.Cells.Item(43, 5).Value = 234.1100006103516
$ws.Cells.Item(43, 6).Value = 236.5299987792969
$ws.Cells.Item(43, 7).Value = 217.9299926757812
$ws.Cells.Item(43, 8).Value = 10664912097
$ws.Cells.Item(43, 9).Value = "AMZN"

